$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.135.12"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "1.787.73"
$ws.Range("E3").Value = "  +0.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.44%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9976"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3961"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3421"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.27"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.92%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07301"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.51%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9979"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.211"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.222"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.25%  "

$ws.Range("D16").Value = "1.785.08"
$ws.Range("E16").Value = "  +0.00%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001070"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.85%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06637"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.79%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9972"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.283"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.24%  "

$ws.Range("D23").Value = "28.133.12"
$ws.Range("E23").Value = "  -0.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.390"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.16%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.20"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.64%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.379"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.43%  "

$ws.Range("D29").Value = "1.976.15"
$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.293"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -8.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "130.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.07%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.074"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.77%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.903"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08785"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.25"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.84%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06275"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02315"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.46%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.205"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.97%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6590"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.31%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2134"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.09%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.507"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.30%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.210"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.18%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9974"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.25%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.38%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6091"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.28%  "

$ws.Range("B47").Value = "PancakeSwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.833"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.67"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.018"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.164"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.47%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07057"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.69%  "
